# Fruta / hortaliza, semanal
# Re-samples the "Ramas de apio" records at Vega Modelo de Temuco onto a
# weekly cadence: for each data row (2-24, row 16 untouched) the
# Fecha / Volumen / Precio minimo / Precio maximo / Precio promedio
# ponderado / Precio $/Kg values are replaced with another row's values
# from the original table (a permutation of the existing observations).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "2021-06-18"
$ws.Range("J2").Value = 55
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = 5000
$ws.Range("P2").Value = 5000
$ws.Range("D3").Value = "2022-03-29"
$ws.Range("J3").Value = 20
$ws.Range("K3").Value = 5000
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = 5000
$ws.Range("P3").Value = 5000
$ws.Range("D4").Value = "2022-04-28"
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = 5000
$ws.Range("P4").Value = 5000
$ws.Range("D5").Value = "2022-07-07"
$ws.Range("J5").Value = 65
$ws.Range("K5").Value = 6000
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = 6000
$ws.Range("P5").Value = 6000
$ws.Range("D6").Value = "2020-12-11"
$ws.Range("J6").Value = 10
$ws.Range("K6").Value = 4000
$ws.Range("L6").Value = 4000
$ws.Range("M6").Value = 4000
$ws.Range("P6").Value = 4000
$ws.Range("D7").Value = "2021-03-25"
$ws.Range("J7").Value = 55
$ws.Range("K7").Value = 4000
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = 4000
$ws.Range("P7").Value = 4000
$ws.Range("D8").Value = "2021-07-13"
$ws.Range("J8").Value = 55
$ws.Range("K8").Value = 6000
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = 6000
$ws.Range("P8").Value = 6000
$ws.Range("D9").Value = "2021-10-28"
$ws.Range("J9").Value = 20
$ws.Range("K9").Value = 4000
$ws.Range("L9").Value = 4000
$ws.Range("M9").Value = 4000
$ws.Range("P9").Value = 4000
$ws.Range("D10").Value = "2022-04-05"
$ws.Range("J10").Value = 85
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 5000
$ws.Range("M10").Value = 5000
$ws.Range("P10").Value = 5000
$ws.Range("D11").Value = "2021-11-08"
$ws.Range("J11").Value = 30
$ws.Range("K11").Value = 4000
$ws.Range("L11").Value = 4000
$ws.Range("M11").Value = 4000
$ws.Range("P11").Value = 4000
$ws.Range("D12").Value = "2021-04-05"
$ws.Range("J12").Value = 35
$ws.Range("K12").Value = 4000
$ws.Range("L12").Value = 4000
$ws.Range("M12").Value = 4000
$ws.Range("P12").Value = 4000
$ws.Range("D13").Value = "2021-04-27"
$ws.Range("J13").Value = 20
$ws.Range("K13").Value = 4000
$ws.Range("L13").Value = 4000
$ws.Range("M13").Value = 4000
$ws.Range("P13").Value = 4000
$ws.Range("D14").Value = "2021-03-04"
$ws.Range("J14").Value = 30
$ws.Range("K14").Value = 4000
$ws.Range("L14").Value = 4000
$ws.Range("M14").Value = 4000
$ws.Range("P14").Value = 4000
$ws.Range("D15").Value = "2021-04-26"
$ws.Range("J15").Value = 50
$ws.Range("K15").Value = 4000
$ws.Range("L15").Value = 4000
$ws.Range("M15").Value = 4000
$ws.Range("P15").Value = 4000
$ws.Range("D17").Value = "2021-04-29"
$ws.Range("J17").Value = 40
$ws.Range("K17").Value = 4000
$ws.Range("L17").Value = 4000
$ws.Range("M17").Value = 4000
$ws.Range("P17").Value = 4000
$ws.Range("D18").Value = "2022-04-29"
$ws.Range("J18").Value = 20
$ws.Range("K18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("M18").Value = 5000
$ws.Range("P18").Value = 5000
$ws.Range("D19").Value = "2022-08-04"
$ws.Range("J19").Value = 25
$ws.Range("K19").Value = 5000
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = 5000
$ws.Range("P19").Value = 5000
$ws.Range("D20").Value = "2021-04-30"
$ws.Range("J20").Value = 20
$ws.Range("K20").Value = 4000
$ws.Range("L20").Value = 4000
$ws.Range("M20").Value = 4000
$ws.Range("P20").Value = 4000
$ws.Range("D21").Value = "2021-04-15"
$ws.Range("J21").Value = 40
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 3000
$ws.Range("M21").Value = 3000
$ws.Range("P21").Value = 3000
$ws.Range("D22").Value = "2021-10-29"
$ws.Range("J22").Value = 40
$ws.Range("K22").Value = 4000
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = 4000
$ws.Range("P22").Value = 4000
$ws.Range("D23").Value = "2022-08-08"
$ws.Range("J23").Value = 40
$ws.Range("K23").Value = 5000
$ws.Range("L23").Value = 5000
$ws.Range("M23").Value = 5000
$ws.Range("P23").Value = 5000
$ws.Range("D24").Value = "2021-11-04"
$ws.Range("J24").Value = 55
$ws.Range("K24").Value = 4000
$ws.Range("L24").Value = 4000
$ws.Range("M24").Value = 4000
$ws.Range("P24").Value = 4000
